$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: SCSE Camp - withdraw 1 slot, register 1 committee slot used
$ws.Cells.Item(2, 8).Value = 89
$ws.Cells.Item(2, 9).Value = 9

# Row 3: ADM Camp - slots unchanged (refresh re-write, same values)
$ws.Cells.Item(3, 8).Value = 90
$ws.Cells.Item(3, 9).Value = 10

# Row 6: EEE Camp - withdraw 1 slot
$ws.Cells.Item(6, 8).Value = 89
$ws.Cells.Item(6, 9).Value = 10

# Row 7: NBS Camp - withdraw 1 slot (down to 0)
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 10

# Row 8: ADM Camp - register 2 more slots (2 -> 4)
$ws.Cells.Item(8, 8).Value = 4
$ws.Cells.Item(8, 9).Value = 10
